# Updated cryptos list on Sat Feb 18 19:49:24 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row,
# and fixes the Polygon/BinanceUSD and ImmutableX/Hedera row ordering.
# Numeric-looking price strings are written with a leading apostrophe so
# Excel keeps them as literal text (matching the workbook's existing
# "24.691.97"-style text prices) instead of silently parsing them into
# floating point numbers and losing formatting (e.g. trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.650.65'
$ws.Range('E2').Value = '  +0.92%  '
$ws.Range('D3').Value = '1.690.32'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  +0.62%  '
$ws.Range('D5').Value = '''315.69'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('D6').Value = '''1.002'
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('D7').Value = '''0.3941'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '''0.4048'
$ws.Range('E8').Value = '  -0.20%  '
$ws.Range('B9').Value = 'BinanceUSD'
$ws.Range('C9').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D9').Value = '''1.003'
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '''1.481'
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('D11').Value = '''52.70'
$ws.Range('E11').Value = '  -3.58%  '
$ws.Range('D12').Value = '''0.08801'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = '''7.213'
$ws.Range('E13').Value = '  -1.68%  '
$ws.Range('D14').Value = '''23.53'
$ws.Range('E14').Value = '  +0.77%  '
$ws.Range('D15').Value = '''8.018'
$ws.Range('E15').Value = '  +7.49%  '
$ws.Range('D16').Value = '''0.00001315'
$ws.Range('E16').Value = '  -0.55%  '
$ws.Range('D17').Value = '1.698.15'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '''99.99'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').Value = '''0.07012'
$ws.Range('E19').Value = '  -0.56%  '
$ws.Range('D20').Value = '''19.44'
$ws.Range('E20').Value = '  +0.09%  '
$ws.Range('D21').Value = '''7.000'
$ws.Range('E21').Value = '  +3.40%  '
$ws.Range('D22').Value = '''1.002'
$ws.Range('E22').Value = '  +0.52%  '
$ws.Range('D23').Value = '''14.22'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').Value = '24.638.76'
$ws.Range('E24').Value = '  +0.90%  '
$ws.Range('D25').Value = '''3.247'
$ws.Range('E25').Value = '  +8.11%  '
$ws.Range('D26').Value = '''2.369'
$ws.Range('E26').Value = '  +2.81%  '
$ws.Range('D27').Value = '''22.72'
$ws.Range('E27').Value = '  +1.12%  '
$ws.Range('D28').Value = '''163.21'
$ws.Range('E28').Value = '  +2.53%  '
$ws.Range('D29').Value = '''135.78'
$ws.Range('E29').Value = '  +1.67%  '
$ws.Range('D30').Value = '''5.193'
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('D31').Value = '''7.552'
$ws.Range('E31').Value = '  +1.12%  '
$ws.Range('D32').Value = '1.885.13'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.08559'
$ws.Range('E33').Value = '  -1.84%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '''1.052'
$ws.Range('E34').Value = '  -3.50%  '
$ws.Range('D35').Value = '''7.133'
$ws.Range('E35').Value = '  -3.36%  '
$ws.Range('D36').Value = '''11.28'
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('D37').Value = '''0.2731'
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('D38').Value = '''1.902'
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('D39').Value = '''14.33'
$ws.Range('E39').Value = '  -2.44%  '
$ws.Range('D40').Value = '''0.09202'
$ws.Range('E40').Value = '  +2.55%  '
$ws.Range('D41').Value = '''0.02716'
$ws.Range('E41').Value = '  -2.72%  '
$ws.Range('D42').Value = '''1.462'
$ws.Range('E42').Value = '  -0.27%  '
$ws.Range('D43').Value = '''0.7586'
$ws.Range('E43').Value = '  -0.28%  '
$ws.Range('D44').Value = '''16.05'
$ws.Range('E44').Value = '  +4.18%  '
$ws.Range('D45').Value = '''0.7121'
$ws.Range('E45').Value = '  -0.72%  '
$ws.Range('D46').Value = '''2.564'
$ws.Range('E46').Value = '  +4.76%  '
$ws.Range('D47').Value = '''4.216'
$ws.Range('E47').Value = '  +1.60%  '
$ws.Range('E48').Value = '  +0.50%  '
$ws.Range('D49').Value = '''1.315'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('D50').Value = '''138.99'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('D51').Value = '''0.07968'
$ws.Range('E51').Value = '  -0.04%  '
